# Updated handling if no traveller is found
# - Fix the stray leading double-space on the "Work/employment..." purpose text
# - Add new Traveller_* columns (doctype, docnum, surname, givename, dob)
#   with sample "no traveller found" fallback data across the three rows
# - Update the active selection to the new traveller name column

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("arrivecan")

# --- Fix column B (TravelPurp_lbl) text: drop the stray leading spaces ---
$ws.Range("B2").Value = "Work/employment or essential reasons  && Cross-border worker"
$ws.Range("B3").Value = "Work/employment or essential reasons  && Cross-border worker"
$ws.Range("B4").Value = "Work/employment or essential reasons  && Cross-border worker"

# --- New headers for the traveller document columns (H1:L1) ---
$ws.Range("H1").Value = "Traveller_doctype"
$ws.Range("I1").Value = "Traveller_docnum"
$ws.Range("J1").Value = "Traveller_surname"
$ws.Range("K1").Value = "Traveller_givename"
$ws.Range("L1").Value = "Traveller_dob"

# --- New traveller data, identical across rows 2-4 ---
# Columns I (docnum) and L (dob) must stay text (same "@" text style the
# existing Traveller_mobile column E uses) so the long id / date-like
# string round-trip verbatim instead of turning into a number or date.
$ws.Range("I2:I4").NumberFormat = "@"
$ws.Range("L2:L4").NumberFormat = "@"

foreach ($r in 2..4) {
    $ws.Range("H$r").Value = "US Permanent Resident Card"
    $ws.Range("I$r").Value = "5647835784358"
    $ws.Range("J$r").Value = "John"
    $ws.Range("K$r").Value = "Eric"
    $ws.Range("L$r").Value = "2020-10-01"
}

# Match the text style (quote-prefixed "@" text format) already used by the
# Traveller_mobile column (E) for the two text-sensitive new columns.
$ws.Range("E2:E4").Copy()
$ws.Range("I2:I4").PasteSpecial(-4122)
$ws.Range("E2:E4").Copy()
$ws.Range("L2:L4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Selection moves to the traveller name column ---
$ws.Range("B2:B4").Select() | Out-Null
